$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: cardholder name / account number change ---
$ws.Range("C2").Value = "Hartmut"
# The card/account number is a 16-digit string that must stay text (it was
# stored as inline text in the source file) -- a leading apostrophe forces
# Excel to keep it as text instead of auto-converting it to a number.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance date ---
$ws.Range("D5").Value = "KONTOSTAND AM 22.03.2024"

# --- Row 6 (was the 02.03./03.03. MITGLIEDSBEITRAG line; now 23.03./24.03. EBAY) ---
$ws.Range("B6").Value = "23.03."
$ws.Range("C6").Value = "24.03."
$ws.Range("D6").Value = "EBAY MKTPLC EU PLRUCY"
$ws.Range("E6").Value = "216,67-"

# --- Row 7 (now ABSCHLAG STROM) ---
$ws.Range("B7").Value = "27.03."
$ws.Range("C7").Value = "28.03."
$ws.Range("D7").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 46825501"
$ws.Range("E7").Value = "84,74-"

# --- Row 8 (now PAYPAL OSQMDW) ---
$ws.Range("B8").Value = "28.03."
$ws.Range("C8").Value = "29.03."
$ws.Range("D8").Value = "PAYPAL OSQMDW"
$ws.Range("E8").Value = "31,53-"

# --- Row 9 (now BEITRAG Allianz SE K-40080526) ---
$ws.Range("B9").Value = "01.04."
$ws.Range("C9").Value = "02.04."
$ws.Range("D9").Value = "BEITRAG Allianz SE K-40080526"
$ws.Range("E9").Value = "52,82-"

# --- Row 10 (now MITGLIEDSBEITRAG ZEUS BODYPOWER) ---
$ws.Range("B10").Value = "02.04."
$ws.Range("C10").Value = "03.04."
$ws.Range("D10").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E10").Value = "24,77-"

# --- Row 11: transaction removed -> cells cleared, E11 style switches from s=17 to s=12 ---
$ws.Range("B11:D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").WrapText = $true
$ws.Range("E11").HorizontalAlignment = -4152

# --- Closing balance date / amount ---
$ws.Range("D12").Value = "KONTOSTAND AM 05.04.2024"
$ws.Range("E12").Value = "410,53-"

# --- Next billing date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 11.04.2024"
